$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" (strikeouts) values per pitching outing.
# Regenerated values (K instead of the old Strike# metric) for rows 2-8 and 11.
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 7
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 4
$ws.Range("G7").Value = 5
$ws.Range("G8").Value = 3
$ws.Range("G11").Value = 4
